$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the category/value pairs between row 1 and row 2
$ws.Range("B1").Value = "D"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "81.32`n"

$ws.Range("B2").Value = "A"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "127.58`n"
